$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Name" row (row 4): fill in the Value column with the resource name
$ws.Range("B4").Value = "StatuthospitalierVs"

# Update the Date row (row 8) value to reflect the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
